# Scheduled runner update: refresh market-board derived price/profit figures
# (currentAveragePrice*, LevePrice*, LeveProfit*) across the per-job Leve sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 33.5
$ws.Range("I6").Value = 32
$ws.Range("K6").Value = 96
$ws.Range("M6").Value = 16

$ws.Range("H58").Value = 642.8
$ws.Range("J58").Value = 799.75
$ws.Range("L58").Value = 2399.25
$ws.Range("N58").Value = -2699.25

$ws.Range("H62").Value = 4621.2
$ws.Range("I62").Value = 5033.3335
$ws.Range("J62").Value = 4003
$ws.Range("K62").Value = 5033.3335
$ws.Range("L62").Value = 4003
$ws.Range("M62").Value = -4409.3335
$ws.Range("N62").Value = -5251

$ws.Range("H65").Value = 4621.2
$ws.Range("I65").Value = 5033.3335
$ws.Range("J65").Value = 4003
$ws.Range("K65").Value = 25166.6675
$ws.Range("L65").Value = 20015
$ws.Range("M65").Value = -22046.6675
$ws.Range("N65").Value = -26255

$ws.Range("H92").Value = 3990.6
$ws.Range("I92").Value = 2488.5
$ws.Range("K92").Value = 2488.5
$ws.Range("M92").Value = -1240.5

$ws.Range("H113").Value = 2168.8
$ws.Range("I113").Value = 2156.2856
$ws.Range("K113").Value = 2156.2856
$ws.Range("M113").Value = 1097.7144

$ws.Range("H136").Value = 89999.75
$ws.Range("J136").Value = 89999.75
$ws.Range("L136").Value = 89999.75
$ws.Range("N136").Value = -100199.75

$ws.Range("H138").Value = 4246.857
$ws.Range("I138").Value = 5874.5
$ws.Range("J138").Value = 3595.8
$ws.Range("K138").Value = 17623.5
$ws.Range("L138").Value = 10787.4
$ws.Range("M138").Value = -12483.5
$ws.Range("N138").Value = -21067.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H15").Value = 4999
$ws.Range("J15").Value = 4999
$ws.Range("L15").Value = 4999
$ws.Range("N15").Value = -5699

$ws.Range("H45").Value = 3945.75
$ws.Range("I45").Value = 2792.6667
$ws.Range("J45").Value = 7405
$ws.Range("K45").Value = 2792.6667
$ws.Range("L45").Value = 7405
$ws.Range("M45").Value = -2415.6667
$ws.Range("N45").Value = -8159

$ws.Range("H122").Value = 2264
$ws.Range("I122").Value = 2264
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 6792
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -4342
$ws.Range("N122").Value = ""

$ws.Range("H132").Value = 8494.5
$ws.Range("I132").Value = 999
$ws.Range("J132").Value = 15990
$ws.Range("K132").Value = 2997
$ws.Range("L132").Value = 47970
$ws.Range("M132").Value = -467
$ws.Range("N132").Value = -53030

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 128.85715
$ws.Range("I22").Value = 142
$ws.Range("J22").Value = 50
$ws.Range("K22").Value = 142
$ws.Range("L22").Value = 50
$ws.Range("M22").Value = 31
$ws.Range("N22").Value = -396

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 233.25
$ws.Range("I7").Value = 299.55554
$ws.Range("J7").Value = 34.333332
$ws.Range("K7").Value = 299.55554
$ws.Range("L7").Value = 34.333332
$ws.Range("M7").Value = -186.55554
$ws.Range("N7").Value = -260.333332

$ws.Range("H22").Value = 275.8
$ws.Range("I22").Value = 319.75
$ws.Range("J22").Value = 100
$ws.Range("K22").Value = 319.75
$ws.Range("L22").Value = 100
$ws.Range("M22").Value = 30.25
$ws.Range("N22").Value = -800

$ws.Range("H31").Value = 15899.8
$ws.Range("J31").Value = 17125
$ws.Range("L31").Value = 17125
$ws.Range("N31").Value = -17715

$ws.Range("H32").Value = 3653.5
$ws.Range("I32").Value = 630.5
$ws.Range("J32").Value = 4661.1665
$ws.Range("K32").Value = 630.5
$ws.Range("L32").Value = 4661.1665
$ws.Range("M32").Value = -314.5
$ws.Range("N32").Value = -5293.1665

$ws.Range("H34").Value = 15899.8
$ws.Range("J34").Value = 17125
$ws.Range("L34").Value = 17125
$ws.Range("N34").Value = -17529

$ws.Range("H132").Value = 11149.5
$ws.Range("I132").Value = 9800
$ws.Range("J132").Value = 12499
$ws.Range("K132").Value = 29400
$ws.Range("L132").Value = 37497
$ws.Range("M132").Value = -26870
$ws.Range("N132").Value = -42557

$ws.Range("H134").Value = 0
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = ""
$ws.Range("N134").Value = ""

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H109").Value = 2383
$ws.Range("I109").Value = 2383
$ws.Range("K109").Value = 7149
$ws.Range("M109").Value = -6109

$ws.Range("H139").Value = 9000
$ws.Range("I139").Value = 9000
$ws.Range("K139").Value = 27000
$ws.Range("M139").Value = -21860

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 52.77778
$ws.Range("I2").Value = 35.714287
$ws.Range("J2").Value = 112.5
$ws.Range("K2").Value = 35.714287
$ws.Range("L2").Value = 112.5
$ws.Range("M2").Value = 77.285713
$ws.Range("N2").Value = -338.5

$ws.Range("H39").Value = 29000
$ws.Range("J39").Value = 29000
$ws.Range("L39").Value = 29000
$ws.Range("N39").Value = -30064

$ws.Range("H97").Value = 2996.3333
$ws.Range("I97").Value = 2990
$ws.Range("J97").Value = 2999.5
$ws.Range("K97").Value = 2990
$ws.Range("L97").Value = 2999.5
$ws.Range("M97").Value = -2494
$ws.Range("N97").Value = -3991.5

$ws.Range("H132").Value = 8400
$ws.Range("I132").Value = 8400
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 25200
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -22670
$ws.Range("N132").Value = ""

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3376.3333
$ws.Range("I122").Value = 3466.3333
$ws.Range("K122").Value = 10398.9999
$ws.Range("M122").Value = -7948.999899999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1851.8636
$ws.Range("J113").Value = 2093.5334
$ws.Range("L113").Value = 6280.600199999999
$ws.Range("N113").Value = -10620.6002

$ws.Range("H132").Value = 10000
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 10000
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 30000
$ws.Range("N132").Value = -35060
$ws.Range("M132").Value = ""
